$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 46 (pushes old rows 46-48 down to 47-49)
$ws.Rows.Item(46).Insert()

# Populate the new row 46 with the new weekly price entry
$ws.Cells.Item(46, 1).Value = 11
$ws.Cells.Item(46, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(46, 3).Value = "Bíobío"
$ws.Cells.Item(46, 4).Value = 44516
$ws.Cells.Item(46, 5).Value = 8
$ws.Cells.Item(46, 6).Value = 100112001
$ws.Cells.Item(46, 7).Value = "Berenjena"
$ws.Cells.Item(46, 8).Value = "Sin especificar"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 350
$ws.Cells.Item(46, 11).Value = 7500
$ws.Cells.Item(46, 12).Value = 8000
$ws.Cells.Item(46, 13).Value = 7714
$ws.Cells.Item(46, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(46, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(46, 16).Value = 129
$ws.Cells.Item(46, 17).Value = 60
$ws.Cells.Item(46, 18).Value = "Hortaliza"
